$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") holds value 59 for every data row (2-79).
# The data was dropped/miscoded during scraping and needs to be corrected to 259.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 79 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 59) {
        $cell.Value = 259
    }
}
